# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H:N) on several rows across all 8 sheets, per the source diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 49
$ws.Range("H49").Value = 1258.5
$ws.Range("I49").Value = 1508
$ws.Range("J49").Value = 1009
$ws.Range("K49").Value = 4524
$ws.Range("L49").Value = 3027
$ws.Range("M49").Value = -4388
$ws.Range("N49").Value = -3299
# Row 87
$ws.Range("H87").Value = 43935.54
$ws.Range("J87").Value = 43935.54
$ws.Range("L87").Value = 43935.54
$ws.Range("N87").Value = -46431.54
# Row 90
$ws.Range("H90").Value = 43935.54
$ws.Range("J90").Value = 43935.54
$ws.Range("L90").Value = 131806.62
$ws.Range("N90").Value = -144286.62
# Row 138
$ws.Range("H138").Value = 3697.6292
$ws.Range("I138").Value = 1825.7097
$ws.Range("J138").Value = 4698.1377
$ws.Range("K138").Value = 5477.1291
$ws.Range("L138").Value = 14094.4131
$ws.Range("M138").Value = -337.1291000000001
$ws.Range("N138").Value = -24374.4131

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9330.569
$ws.Range("I32").Value = 8550.433999999999
$ws.Range("J32").Value = 17600
$ws.Range("K32").Value = 8550.433999999999
$ws.Range("L32").Value = 17600
$ws.Range("M32").Value = -8263.433999999999
$ws.Range("N32").Value = -18174
# Row 64
$ws.Range("H64").Value = 19882
$ws.Range("I64").Value = 19882
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 19882
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -19634
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 19882
$ws.Range("I67").Value = 19882
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 19882
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -19024
$ws.Range("N67").ClearContents()
# Row 74
$ws.Range("H74").Value = 115635.25
$ws.Range("I74").Value = 124037.97
$ws.Range("J74").Value = 34409
$ws.Range("K74").Value = 124037.97
$ws.Range("L74").Value = 34409
$ws.Range("M74").Value = -123163.97
$ws.Range("N74").Value = -36157
# Row 77
$ws.Range("H77").Value = 115635.25
$ws.Range("I77").Value = 124037.97
$ws.Range("J77").Value = 34409
$ws.Range("K77").Value = 620189.85
$ws.Range("L77").Value = 172045
$ws.Range("M77").Value = -615821.85
$ws.Range("N77").Value = -180781
# Row 97
$ws.Range("H97").Value = 1251.9667
$ws.Range("I97").Value = 964.7143
$ws.Range("J97").Value = 1922.2222
$ws.Range("K97").Value = 964.7143
$ws.Range("L97").Value = 1922.2222
$ws.Range("M97").Value = -468.7143
$ws.Range("N97").Value = -2914.2222
# Row 132
$ws.Range("H132").Value = 2828.325
$ws.Range("I132").Value = 1663.1904
$ws.Range("K132").Value = 4989.5712
$ws.Range("M132").Value = -2459.5712
# Row 140
$ws.Range("H140").Value = 36705.848
$ws.Range("J140").Value = 37854.184
$ws.Range("L140").Value = 37854.184
$ws.Range("N140").Value = -48214.184

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 33454
$ws.Range("I62").Value = 30000
$ws.Range("J62").Value = 35181
$ws.Range("K62").Value = 30000
$ws.Range("L62").Value = 35181
$ws.Range("M62").Value = -29314
$ws.Range("N62").Value = -36553
# Row 65
$ws.Range("H65").Value = 33454
$ws.Range("I65").Value = 30000
$ws.Range("J65").Value = 35181
$ws.Range("K65").Value = 90000
$ws.Range("L65").Value = 105543
$ws.Range("M65").Value = -86568
$ws.Range("N65").Value = -112407
# Row 94
$ws.Range("H94").Value = 2246.8125
$ws.Range("I94").Value = 1803.9
$ws.Range("K94").Value = 1803.9
$ws.Range("M94").Value = -1352.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3793.2
$ws.Range("I31").Value = 2077.2554
$ws.Range("J31").Value = 8273.723
$ws.Range("K31").Value = 2077.2554
$ws.Range("L31").Value = 8273.723
$ws.Range("M31").Value = -1782.2554
$ws.Range("N31").Value = -8863.723
# Row 32
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 10000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9684
$ws.Range("N32").ClearContents()
# Row 34
$ws.Range("H34").Value = 3793.2
$ws.Range("I34").Value = 2077.2554
$ws.Range("J34").Value = 8273.723
$ws.Range("K34").Value = 2077.2554
$ws.Range("L34").Value = 8273.723
$ws.Range("M34").Value = -1875.2554
$ws.Range("N34").Value = -8677.723

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 34.696968
$ws.Range("I2").Value = 33.857143
$ws.Range("J2").Value = 34.923077
$ws.Range("K2").Value = 203.142858
$ws.Range("L2").Value = 209.538462
$ws.Range("M2").Value = -90.14285799999999
$ws.Range("N2").Value = -435.538462
# Row 26
$ws.Range("H26").Value = 512.125
$ws.Range("I26").Value = 75.625
$ws.Range("J26").Value = 948.625
$ws.Range("K26").Value = 226.875
$ws.Range("L26").Value = 2845.875
$ws.Range("M26").Value = 61.125
$ws.Range("N26").Value = -3421.875
# Row 32
$ws.Range("H32").Value = 1996.24
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 2082.8696
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 6248.6088
$ws.Range("M32").Value = -2717
$ws.Range("N32").Value = -6814.6088
# Row 38
$ws.Range("H38").Value = 89.09090999999999
$ws.Range("I38").Value = 38.333332
$ws.Range("J38").Value = 150
$ws.Range("K38").Value = 114.999996
$ws.Range("L38").Value = 450
$ws.Range("M38").Value = 232.000004
$ws.Range("N38").Value = -1144
# Row 39
$ws.Range("H39").Value = 9948
$ws.Range("J39").Value = 9948
$ws.Range("L39").Value = 29844
$ws.Range("N39").Value = -30432
# Row 46
$ws.Range("H46").Value = 2737.5
$ws.Range("J46").Value = 3085.7144
$ws.Range("L46").Value = 9257.143199999999
$ws.Range("N46").Value = -9439.143199999999
# Row 51
$ws.Range("H51").Value = 1999.5
$ws.Range("J51").Value = 1999.5
$ws.Range("L51").Value = 5998.5
$ws.Range("N51").Value = -6918.5
# Row 126
$ws.Range("H126").Value = 1605.2354
$ws.Range("I126").Value = 1090.75
$ws.Range("J126").Value = 2840
$ws.Range("K126").Value = 3272.25
$ws.Range("L126").Value = 8520
$ws.Range("M126").Value = 1667.75
$ws.Range("N126").Value = -18400
# Row 131
$ws.Range("H131").Value = 499.61615
$ws.Range("I131").Value = 275.2258
$ws.Range("J131").Value = 875.62164
$ws.Range("K131").Value = 825.6774
$ws.Range("L131").Value = 2626.86492
$ws.Range("M131").Value = 4214.3226
$ws.Range("N131").Value = -12706.86492

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 6499.7856
$ws.Range("I80").Value = 14524.5
$ws.Range("K80").Value = 14524.5
$ws.Range("M80").Value = -13526.5
# Row 83
$ws.Range("H83").Value = 6499.7856
$ws.Range("I83").Value = 14524.5
$ws.Range("K83").Value = 72622.5
$ws.Range("M83").Value = -67630.5
# Row 132
$ws.Range("H132").Value = 26394.523
$ws.Range("I132").Value = 68240.47
$ws.Range("J132").Value = 3146.7778
$ws.Range("K132").Value = 204721.41
$ws.Range("L132").Value = 9440.3334
$ws.Range("M132").Value = -202191.41
$ws.Range("N132").Value = -14500.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1217.5
$ws.Range("I16").Value = 1105.7142
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1105.7142
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -935.7141999999999
$ws.Range("N16").Value = -2340
# Row 132
$ws.Range("H132").Value = 2993.611
$ws.Range("I132").Value = 2302.4285
$ws.Range("K132").Value = 6907.2855
$ws.Range("M132").Value = -4377.2855
# Row 136
$ws.Range("H136").Value = 4271.703
$ws.Range("I136").Value = 2696.375
$ws.Range("J136").Value = 6897.25
$ws.Range("K136").Value = 8089.125
$ws.Range("L136").Value = 20691.75
$ws.Range("M136").Value = -5539.125
$ws.Range("N136").Value = -25791.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 8500
# Row 129
$ws.Range("H129").Value = 52450
$ws.Range("J129").Value = 52450
$ws.Range("L129").Value = 52450
$ws.Range("N129").Value = -62450
# Row 132
$ws.Range("H132").Value = 2783.3901
$ws.Range("I132").Value = 2416.9033
$ws.Range("J132").Value = 3919.5
$ws.Range("K132").Value = 7250.7099
$ws.Range("L132").Value = 11758.5
$ws.Range("M132").Value = -4720.7099
$ws.Range("N132").Value = -16818.5
# Row 136
$ws.Range("H136").Value = 3869.3594
$ws.Range("I136").Value = 1628.2979
$ws.Range("J136").Value = 10065.235
$ws.Range("K136").Value = 4884.893700000001
$ws.Range("L136").Value = 30195.705
$ws.Range("M136").Value = -2334.893700000001
$ws.Range("N136").Value = -35295.705
